$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 (SVM, All Features PCA): update Test Accuracy
$ws.Range("F6").Value = 0.8582941571524513

# Row 7 (SVM, Individual Features PCA): update hyperparameters, train accuracy, mean CV accuracy
$ws.Range("C7").Value = "C: 10, class_weight: balanced, degree: 2, gamma: 0.001, kernel: rbf"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.8062611712461838

# Row 9 (XGBoost, All Features): update hyperparameters, mean CV accuracy
$ws.Range("C9").Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 1, n_estimators: 300"
$ws.Range("E9").Value = 0.9085460727171804
